# Update the "K" column (column G) of the save_data sheet with the
# newly regenerated strikeout values (K instead of the old Strike# calc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 0
    4  = 2
    6  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 2
    12 = 0
    13 = 0
    14 = 2
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 1
    20 = 1
    21 = 0
    22 = 1
    23 = 2
    24 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 0
    32 = 2
    33 = 0
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 1
    39 = 2
    40 = 2
    41 = 3
    42 = 2
    43 = 0
    45 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
